$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "EObjectImpl@3c72031c" "EObjectImpl@1bd5577e"
Replace-Text "M2DocEvaluator.java:543)" "M2DocEvaluator.java:555)"
Replace-Text "M2DocEvaluator.java:1084)" "M2DocEvaluator.java:1096)"
Replace-Text "M2DocEvaluator.java:1300)" "M2DocEvaluator.java:1305)"
Replace-Text "M2DocEvaluator.java:278)" "M2DocEvaluator.java:283)"
Replace-Text "M2DocEvaluator.java:267)" "M2DocEvaluator.java:272)"
Replace-Text "AbstractTemplatesTestSuite.java:476)" "AbstractTemplatesTestSuite.java:480)"
Replace-Text "AbstractTemplatesTestSuite.java:385)" "AbstractTemplatesTestSuite.java:389)"
Replace-Text "GeneratedMethodAccessor111" "GeneratedMethodAccessor107"
